# "Add files via upload" — the workbook's single sheet (Sheet1, a wide
# A:AL intake-form export) is replaced by a two-sheet "Client" export:
#   - Client0000: a small 3-column (Address / Employment / Education) sheet
#   - Client0001: the bulk of the original columns, reordered/renamed, with
#                 the two original data rows folded into a single row.

function New-RowArray {
    param($Values)
    $n = $Values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    return $arr
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: reuse the existing sheet, rename + replace its content ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Client0000"
$ws1.Cells.Clear()

$sheet1Headers = New-RowArray @("Address", "Employment", "Education")
$ws1.Range("A1").Resize(1, 3).Value = $sheet1Headers

# Row 2 stays empty (no values) but is still part of the sheet's used range
# (dimension A1:C2) — touch it with a no-op border tweak (xlLineStyleNone)
# so it registers without writing any content or cell styling.
$ws1.Range("A2:C2").Borders.LineStyle = -4142

# --- Sheet 2: brand-new sheet, inserted right after Client0000 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Client0001"

$headers = @(
    "A-Number", "SSN", "USCIS Account Number", "Last Name", "First Name",
    "Middle Name", "Other Names", "Current Address", "Mailing Address",
    "Gender", "Marital Status", "Date of Birth", "City of Birth",
    "Country of Birth", "Nationality", "Religion", "Last Leave country",
    "List entry to US: Date", "List entry to US: i94",
    "List entry to US: Place", "List entry to US: Visa Category",
    "List entry to US: Status Expire", "Passport Number",
    "Passport Expiration Date", "Address", "Employment", "Education"
)
$sheet2Headers = New-RowArray $headers
$ws2.Range("A1").Resize(1, $headers.Length).Value = $sheet2Headers

# Single data row — the merge of the original sheet's two data rows
# (row 2's blanks filled in by row 3's Yangon/Myanmar/Burmese/Buddhism values).
$row2 = @(
    "N/A", "N/A", "N/A", "", "", "", "", "", "", "", "", "",
    "Yangon", "Myanmar", "Burmese", "Buddhism", "", "", "", "", "", "",
    "", "", "[]", "[]", "[]"
)
$sheet2Row2 = New-RowArray $row2
$ws2.Range("A2").Resize(1, $row2.Length).Value = $sheet2Row2
